$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue $ws.Range("D2") "63.275.44"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3: Ethereum
Set-TextValue $ws.Range("D3") "2.554.77"
$ws.Range("E3").Value = "  +3.46%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5: BNB
Set-TextValue $ws.Range("D5") "568.77"
$ws.Range("E5").Value = "  +0.78%  "

# Row 6: Solana
Set-TextValue $ws.Range("D6") "147.29"
$ws.Range("E6").Value = "  +3.39%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.08%  "

# Row 8: XRP
$ws.Range("E8").Value = "  -0.08%  "

# Row 9: LidoStakedEther
Set-TextValue $ws.Range("D9") "2.553.31"
$ws.Range("E9").Value = "  +3.47%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +0.80%  "

# Row 11: Toncoin
Set-TextValue $ws.Range("D11") "5.59"
$ws.Range("E11").Value = "  -1.76%  "

# Row 12: TRON
$ws.Range("E12").Value = "  +0.46%  "

# Row 13: Cardano
Set-TextValue $ws.Range("D13") "0.353"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14: Avalanche
Set-TextValue $ws.Range("D14") "27.61"
$ws.Range("E14").Value = "  +3.70%  "

# Row 15: WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.008.11"
$ws.Range("E15").Value = "  +3.38%  "

# Row 16: WrappedBTC
Set-TextValue $ws.Range("D16") "63.175.73"
$ws.Range("E16").Value = "  +0.60%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  +2.27%  "

# Row 18: WrappedEther
Set-TextValue $ws.Range("D18") "2.569.97"
$ws.Range("E18").Value = "  +4.42%  "

# Row 19: Chainlink
$ws.Range("E19").Value = "  +1.98%  "

# Row 20: BitcoinCash
Set-TextValue $ws.Range("D20") "336.22"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21: Polkadot
Set-TextValue $ws.Range("D21") "4.32"
$ws.Range("E21").Value = "  +1.54%  "

# Row 22: Uniswap
Set-TextValue $ws.Range("D22") "6.79"
$ws.Range("E22").Value = "  -0.26%  "

# Row 23: Dai
$ws.Range("E23").Value = "  -0.11%  "

# Row 24: Litecoin
Set-TextValue $ws.Range("D24") "65.31"
$ws.Range("E24").Value = "  -0.45%  "

# Row 25: Fetch.AI
Set-TextValue $ws.Range("D25") "1.64"
$ws.Range("E25").Value = "  +9.23%  "

# Row 26: Kaspa
$ws.Range("E26").Value = "  -1.90%  "

# Row 27: SuiNetwork
$ws.Range("E27").Value = "  +7.99%  "

# Row 28: InternetComputer(DFINITY)
Set-TextValue $ws.Range("D28") "8.49"
$ws.Range("E28").Value = "  +5.49%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  -0.03%  "

# Row 30: Aptos
Set-TextValue $ws.Range("D30") "7.40"
$ws.Range("E30").Value = "  +8.20%  "

# Row 31: PEPE
Set-TextValue $ws.Range("D31") "0.0₃0824"
$ws.Range("E31").Value = "  +3.37%  "

# Row 32: PancakeSwap
$ws.Range("E32").Value = "  +1.24%  "

# Row 33: Monero
Set-TextValue $ws.Range("D33") "176.72"
$ws.Range("E33").Value = "  +0.46%  "

# Row 34: ImmutableX
Set-TextValue $ws.Range("D34") "1.58"
$ws.Range("E34").Value = "  +4.87%  "

# Row 35: Bittensor
Set-TextValue $ws.Range("D35") "412.25"
$ws.Range("E35").Value = "  +11.95%  "

# Row 36: PolygonEcosystemToken
Set-TextValue $ws.Range("D36") "0.401"
$ws.Range("E36").Value = "  +1.47%  "

# Row 37: EthereumClassic
Set-TextValue $ws.Range("D37") "19.00"
$ws.Range("E37").Value = "  +1.18%  "

# Row 38: NEARProtocol
$ws.Range("E38").Value = "  +0.63%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  +4.19%  "

# Row 41: FirstDigitalUSD
$ws.Range("E41").Value = "  -0.06%  "

# Row 42: OKB
Set-TextValue $ws.Range("D42") "39.34"
$ws.Range("E42").Value = "  -3.09%  "

# Row 43: Aave
Set-TextValue $ws.Range("D43") "153.46"
$ws.Range("E43").Value = "  +2.54%  "

# Row 44: Filecoin
Set-TextValue $ws.Range("D44") "3.79"
$ws.Range("E44").Value = "  +2.33%  "

# Row 45: InjectiveProtocol
Set-TextValue $ws.Range("D45") "21.17"
$ws.Range("E45").Value = "  +3.47%  "

# Row 46: Mantle
Set-TextValue $ws.Range("D46") "0.607"
$ws.Range("E46").Value = "  +1.63%  "

# Row 47: Stellar
Set-TextValue $ws.Range("D47") "0.0963"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48: Hedera
$ws.Range("E48").Value = "  +1.91%  "

# Row 49: VeChain
Set-TextValue $ws.Range("D49") "0.0240"
$ws.Range("E49").Value = "  +6.33%  "

# Row 50: EnergySwap
Set-TextValue $ws.Range("D50") "18.43"
$ws.Range("E50").Value = "  +2.87%  "

# Row 51: dogwifhat
$ws.Range("E51").Value = "  +2.49%  "
